$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update theta_se row (row 4): B4..G4
$ws.Range("B4").Value = "(0.02)"
$ws.Range("C4").Value = "(0.14)"
$ws.Range("D4").Value = "(0.12)"
$ws.Range("E4").Value = "(0.18)"
$ws.Range("F4").Value = "(0.22)"
$ws.Range("G4").Value = "(0.64)"

# Update lambda_se row (row 6): B6..G6
$ws.Range("B6").Value = "(0.08)"
$ws.Range("C6").Value = "(0.07)"
$ws.Range("D6").Value = "(0.0)"
$ws.Range("E6").Value = "(0.11)"
$ws.Range("F6").Value = "(0.06)"
$ws.Range("G6").Value = "(0.27)"
